$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.269.50"
$ws.Range("E2").Value = "  -5.28%  "
$ws.Range("D3").Value = "2.903.57"
$ws.Range("E3").Value = "  -2.59%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'547.48"
$ws.Range("E5").Value = "  -2.53%  "
$ws.Range("D6").Value = "'124.73"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "2.900.31"
$ws.Range("E8").Value = "  -2.58%  "
$ws.Range("D9").Value = "'0.501"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("E10").Value = "  -8.16%  "
$ws.Range("D11").Value = "'4.66"
$ws.Range("E11").Value = "  -9.30%  "
$ws.Range("D12").Value = "'0.433"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").Value = "'0.0000209"
$ws.Range("E13").Value = "  -6.08%  "
$ws.Range("D14").Value = "'32.11"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "3.385.57"
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("D17").Value = "2.901.02"
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("D18").Value = "'6.48"
$ws.Range("E18").Value = "  +5.73%  "
$ws.Range("D19").Value = "57.265.63"
$ws.Range("E19").Value = "  -5.40%  "
$ws.Range("D20").Value = "'404.15"
$ws.Range("E20").Value = "  -6.15%  "
$ws.Range("D21").Value = "'12.76"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("D22").Value = "'0.669"
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("D23").Value = "'6.79"
$ws.Range("E23").Value = "  -4.63%  "
$ws.Range("D24").Value = "'12.59"
$ws.Range("E24").Value = "  -2.45%  "
$ws.Range("D25").Value = "'77.39"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").Value = "'7.18"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").Value = "'1.91"
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'24.60"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'5.92"
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("D33").Value = "'0.0980"
$ws.Range("E33").Value = "  +4.54%  "
$ws.Range("D34").Value = "'0.912"
$ws.Range("E34").Value = "  -3.83%  "
$ws.Range("D35").Value = "'5.41"
$ws.Range("E35").Value = "  -1.84%  "
$ws.Range("D36").Value = "'1.99"
$ws.Range("E36").Value = "  -12.10%  "
$ws.Range("D37").Value = "'48.11"
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("D38").Value = "'8.20"
$ws.Range("E38").Value = "  +5.69%  "
$ws.Range("D39").Value = "0.0₃0624"
$ws.Range("E39").Value = "  -6.03%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").Value = "'0.0335"
$ws.Range("E41").Value = "  -6.33%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.623.18"
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.41"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").Value = "'360.53"
$ws.Range("E44").Value = "  -3.62%  "
$ws.Range("D46").Value = "'119.34"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.107"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("B48").Value = "TheGraph"
$ws.Range("C48").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D48").Value = "'0.227"
$ws.Range("E48").Value = "  -3.25%  "
$ws.Range("D49").Value = "'1.92"
$ws.Range("E49").Value = "  -2.05%  "
$ws.Range("D50").Value = "'22.45"
$ws.Range("E50").Value = "  -3.74%  "
$ws.Range("D51").Value = "'1.94"
$ws.Range("E51").Value = "  -3.84%  "
